$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. all_products: fill in the (previously blank) sub_category /
#    sub_category_fk columns (H, I) for the product rows.
# ---------------------------------------------------------------------------
$allProducts = $wb.Worksheets.Item("all_products")

$subCategoryData = @{
    2 = @("Sharing", 5)
    3 = @("Sharing", 5)
    4 = @("Nuts", 14)
    5 = @("Sharing", 5)
    6 = @("Dips", 13)
    7 = @("MP", 7)
}

foreach ($row in $subCategoryData.Keys) {
    $vals = $subCategoryData[$row]
    $allProducts.Cells.Item($row, 8).Value = $vals[0]
    $allProducts.Cells.Item($row, 9).Value = $vals[1]
}

# ---------------------------------------------------------------------------
# 2. scif: fill in the (previously blank) VLOOKUP formulas that pull the
#    sub_category / sub_category_fk brand-of-segment columns (P, Q) through
#    from all_products.
# ---------------------------------------------------------------------------
$scif = $wb.Worksheets.Item("scif")

for ($row = 2; $row -le 12; $row++) {
    $scif.Cells.Item($row, 16).Formula = "=VLOOKUP(`$B$row, all_products!`$A`$2:`$O`$14,8, 0)"
    $scif.Cells.Item($row, 17).Formula = "=VLOOKUP(`$B$row, all_products!`$A`$2:`$O`$14,9, 0)"
}

# ---------------------------------------------------------------------------
# 3. matches: the autofilter range shrinks from A1:Q52 to A1:P52 (the filter
#    database no longer spans the now-unused last column).
# ---------------------------------------------------------------------------
$matches = $wb.Worksheets.Item("matches")
$matches.AutoFilterMode = $false
$matches.Range("A1:P52").AutoFilter() | Out-Null

# ---------------------------------------------------------------------------
# 4. Cosmetic view-state tweaks that came along with the edit.
# ---------------------------------------------------------------------------
$matches.Application.Goto($matches.Range("A23"))
$matches.Range("A44").Select() | Out-Null

$allProducts.Application.Goto($allProducts.Range("A1"))
$allProducts.Range("H15").Select() | Out-Null

$scif.Application.Goto($scif.Range("A1"))
$scif.Range("A1").Select() | Out-Null
